# Rework the "UniversalParams" sheet: drop the old WorkingDirectory /
# CellFileName / *Threshold columns, reorder the surviving parameters, and
# add three new ones (LargerSpatialBin, trials_corrTemplate, ds_factor).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Wipe all existing columns (data + column-width formatting) so we can lay
# the new, reordered table down on a clean sheet.
$ws.Columns("A:Z").Delete()

$headers = @("TimeBin","SpatialBin","LargerSpatialBin","TrackStart","TrackEnd","SpeedCutoff","SmoothSigmaFR","trials_corrTemplate","ds_factor")
$values  = @(0.02, 2, 4, 0, 400, 2, 10, 50, 1)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
    $ws.Cells.Item(2, $i + 1).Value = $values[$i]
}

# Column widths for C:H (A/B and I keep the sheet's default width).
$ws.Columns("C:C").ColumnWidth = 12.330729166666666   # -> ~13.1640625
$ws.Columns("D:D").ColumnWidth = 8.166666666666666    # -> 9
$ws.Columns("E:E").ColumnWidth = 13.166666666666666   # -> 14
$ws.Columns("F:F").ColumnWidth = 12.498697916666666   # -> ~13.33203125
$ws.Columns("G:G").ColumnWidth = 14.0                 # -> ~14.83203125
$ws.Columns("H:H").ColumnWidth = 21.166666666666668   # -> 22

# Give the new "trials_corrTemplate" header (H1) its own look: 10pt Helvetica.
$ws.Range("H1").Font.Name = "Helvetica"
$ws.Range("H1").Font.Size = 10

# Match the saved selection/active cell from the edit.
$ws.Range("I1").Select() | Out-Null

Write-Output "edit applied"
